# Added attendance for Monique presentation
# (seminar on "Safety for Whom? ..." row, 2024 - Fall sheet, row 11)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024 - Fall")

# Undergrad (G) / grad (H) attendance counts for the Nov. 6th seminar
# (the day after the 2024 election) — previously recorded with 0
# attendees, now filled in with the real headcount.
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 3

# Update the "room/occasion" note for that session to reflect the date.
$ws.Range("J11").Value = "JMC, Day after 2024 Election"

# Leave the cursor where the author left it when they saved.
$ws.Range("H13").Select()
